# Assignment 14 "User Story Specs.xlsx" update:
# Fill in the previously-blank "Should Have" (column H) cells for the two
# user stories in rows 8 and 9 with the "should not have" field data, and
# move the active selection to the last-edited cell (H9), matching the
# author's edit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 / column H ("Should Have" for "I, As a User of this application...")
$ws.Range("H8").Value = "1. Broken Links within the drop down that do not cause reaction on the page.`n2. Links that are not relevant to the groups within the drop downs."

# Row 9 / column H ("Should Have" for the tabular-display user story)
$ws.Range("H9").Value = "1. Information other than what is required in the table. Simplify what the user sees.`n2. Differing color scheme to conflict with the rest of the page."

# Reflect the scroll position / selection the author ended the session on.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("H9").Select()

# Best-effort: the author's window was resized (maximized) while editing;
# not all COM engines persist this back to bookViews, but set it anyway.
try {
    $win.WindowState = -4137
    $win.Width = 30520
    $win.Height = 17560
} catch {
}
